$wb = $excel.ActiveWorkbook

# --- Sheet1: drop the extra tail rows (A45:A87, the leftover index column)
# that weren't meant to be part of the dataset, shrinking the used range
# back down to A1:N44.
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Rows("45:87").Delete()

# Sheet1 becomes the active/selected sheet, scrolled down near the bottom
# of the (now shorter) data with F65 selected.
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 45
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("F65").Select()

# Note: Sheet3 simply stops being the tab-selected sheet as a side effect of
# Sheet1 becoming active above (only one sheet can be "tabSelected"); its own
# selection (A2:N44) is left untouched, so no further action is needed here.
